$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 used to hold the "no" anchor under the negative-word table; that
# anchor was dropped, so clear the stale A7:H7 cells (and their style).
$ws.Range("A7:H7").Clear()

# Apply the header-row style (bold, bordered, centered) to the new J
# column rows (28-30) before filling them in, matching the style used
# by the existing rows in that table.
$ws.Range("J2").Copy() | Out-Null
$ws.Range("J28:J30").PasteSpecial(-4122) | Out-Null

# Populate the cells with the refreshed anchor-score results.
$ws.Range("A1").Value = "negative"
$ws.Range("J1").Value = "positive"
$ws.Range("A2").Value = "name"
$ws.Range("B2").Value = "anchor score"
$ws.Range("C2").Value = "type occurences"
$ws.Range("D2").Value = "total occurences"
$ws.Range("E2").Value = "+%"
$ws.Range("F2").Value = "-%"
$ws.Range("G2").Value = "both"
$ws.Range("H2").Value = "normal"
$ws.Range("J2").Value = "name"
$ws.Range("K2").Value = "anchor score"
$ws.Range("L2").Value = "type occurences"
$ws.Range("M2").Value = "total occurences"
$ws.Range("N2").Value = "+%"
$ws.Range("O2").Value = "-%"
$ws.Range("P2").Value = "both"
$ws.Range("Q2").Value = "normal"
$ws.Range("A3").Value = "crude"
$ws.Range("B3").Value = 0.8529411764705882
$ws.Range("C3").Value = 29
$ws.Range("D3").Value = 29
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = 5
$ws.Range("J3").Value = "best"
$ws.Range("K3").Value = 0.9491525423728814
$ws.Range("L3").Value = 56
$ws.Range("M3").Value = 56
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = $false
$ws.Range("Q3").Value = 3
$ws.Range("A4").Value = "crisis"
$ws.Range("B4").Value = 0.6095890410958904
$ws.Range("C4").Value = 178
$ws.Range("D4").Value = 178
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = 114
$ws.Range("J4").Value = "interesting"
$ws.Range("K4").Value = 0.9090909090909091
$ws.Range("L4").Value = 30
$ws.Range("M4").Value = 30
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = $false
$ws.Range("Q4").Value = 3
$ws.Range("A5").Value = "panic"
$ws.Range("B5").Value = 0.1821705426356589
$ws.Range("C5").Value = 94
$ws.Range("D5").Value = 94
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = 422
$ws.Range("J5").Value = "love"
$ws.Range("K5").Value = 0.8913043478260869
$ws.Range("L5").Value = 41
$ws.Range("M5").Value = 41
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = $false
$ws.Range("Q5").Value = 5
$ws.Range("A6").Value = "sc"
$ws.Range("B6").Value = 0.1746031746031746
$ws.Range("C6").Value = 33
$ws.Range("D6").Value = 33
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = 156
$ws.Range("J6").Value = "great"
$ws.Range("K6").Value = 0.8660714285714286
$ws.Range("L6").Value = 97
$ws.Range("M6").Value = 97
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = $false
$ws.Range("Q6").Value = 15
$ws.Range("J7").Value = "thanks"
$ws.Range("K7").Value = 0.8292682926829268
$ws.Range("L7").Value = 68
$ws.Range("M7").Value = 68
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = $false
$ws.Range("Q7").Value = 14
$ws.Range("J8").Value = "thank"
$ws.Range("K8").Value = 0.796875
$ws.Range("L8").Value = 102
$ws.Range("M8").Value = 102
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = $false
$ws.Range("Q8").Value = 26
$ws.Range("J9").Value = "positive"
$ws.Range("K9").Value = 0.7931034482758621
$ws.Range("L9").Value = 46
$ws.Range("M9").Value = 46
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = $false
$ws.Range("Q9").Value = 12
$ws.Range("J10").Value = "special"
$ws.Range("K10").Value = 0.7777777777777778
$ws.Range("L10").Value = 28
$ws.Range("M10").Value = 28
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = $false
$ws.Range("Q10").Value = 8
$ws.Range("J11").Value = "won"
$ws.Range("K11").Value = 0.7692307692307693
$ws.Range("L11").Value = 30
$ws.Range("M11").Value = 30
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = $false
$ws.Range("Q11").Value = 9
$ws.Range("J12").Value = "free"
$ws.Range("K12").Value = 0.7416666666666667
$ws.Range("L12").Value = 89
$ws.Range("M12").Value = 89
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = $false
$ws.Range("Q12").Value = 31
$ws.Range("J13").Value = "safe"
$ws.Range("K13").Value = 0.7112676056338029
$ws.Range("L13").Value = 101
$ws.Range("M13").Value = 101
$ws.Range("N13").Value = 1
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = $false
$ws.Range("Q13").Value = 41
$ws.Range("J14").Value = "safety"
$ws.Range("K14").Value = 0.7058823529411765
$ws.Range("L14").Value = 36
$ws.Range("M14").Value = 36
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = $false
$ws.Range("Q14").Value = 15
$ws.Range("J15").Value = "confidence"
$ws.Range("K15").Value = 0.6944444444444444
$ws.Range("L15").Value = 25
$ws.Range("M15").Value = 25
$ws.Range("N15").Value = 1
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = $false
$ws.Range("Q15").Value = 11
$ws.Range("J16").Value = "good"
$ws.Range("K16").Value = 0.69375
$ws.Range("L16").Value = 111
$ws.Range("M16").Value = 111
$ws.Range("N16").Value = 1
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = $false
$ws.Range("Q16").Value = 49
$ws.Range("J17").Value = "support"
$ws.Range("K17").Value = 0.6792452830188679
$ws.Range("L17").Value = 72
$ws.Range("M17").Value = 72
$ws.Range("N17").Value = 1
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = $false
$ws.Range("Q17").Value = 34
$ws.Range("J18").Value = "better"
$ws.Range("K18").Value = 0.6190476190476191
$ws.Range("L18").Value = 39
$ws.Range("M18").Value = 39
$ws.Range("N18").Value = 1
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = $false
$ws.Range("Q18").Value = 24
$ws.Range("J19").Value = "well"
$ws.Range("K19").Value = 0.6170212765957447
$ws.Range("L19").Value = 58
$ws.Range("M19").Value = 58
$ws.Range("N19").Value = 1
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = $false
$ws.Range("Q19").Value = 36
$ws.Range("J20").Value = "heroes"
$ws.Range("K20").Value = 0.6170212765957447
$ws.Range("L20").Value = 29
$ws.Range("M20").Value = 29
$ws.Range("N20").Value = 1
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = $false
$ws.Range("Q20").Value = 18
$ws.Range("J21").Value = "fresh"
$ws.Range("K21").Value = 0.5625
$ws.Range("L21").Value = 27
$ws.Range("M21").Value = 27
$ws.Range("N21").Value = 1
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = $false
$ws.Range("Q21").Value = 21
$ws.Range("J22").Value = "hand"
$ws.Range("K22").Value = 0.5143603133159269
$ws.Range("L22").Value = 197
$ws.Range("M22").Value = 197
$ws.Range("N22").Value = 1
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = $false
$ws.Range("Q22").Value = 186
$ws.Range("J23").Value = "like"
$ws.Range("K23").Value = 0.45
$ws.Range("L23").Value = 153
$ws.Range("M23").Value = 153
$ws.Range("N23").Value = 1
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = $false
$ws.Range("Q23").Value = 187
$ws.Range("J24").Value = "help"
$ws.Range("K24").Value = 0.4406779661016949
$ws.Range("L24").Value = 130
$ws.Range("M24").Value = 130
$ws.Range("N24").Value = 1
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = $false
$ws.Range("Q24").Value = 165
$ws.Range("J25").Value = "care"
$ws.Range("K25").Value = 0.4382022471910113
$ws.Range("L25").Value = 39
$ws.Range("M25").Value = 39
$ws.Range("N25").Value = 1
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = $false
$ws.Range("Q25").Value = 50
$ws.Range("J26").Value = "protect"
$ws.Range("K26").Value = 0.3972602739726027
$ws.Range("L26").Value = 29
$ws.Range("M26").Value = 29
$ws.Range("N26").Value = 1
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = $false
$ws.Range("Q26").Value = 44
$ws.Range("J27").Value = "please"
$ws.Range("K27").Value = 0.3221757322175732
$ws.Range("L27").Value = 77
$ws.Range("M27").Value = 77
$ws.Range("N27").Value = 1
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = $false
$ws.Range("Q27").Value = 162
$ws.Range("J28").Value = "you"
$ws.Range("K28").Value = 0.02083333333333333
$ws.Range("L28").Value = 25
$ws.Range("M28").Value = 25
$ws.Range("N28").Value = 1
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = $false
$ws.Range("Q28").Value = 1175
$ws.Range("J29").Value = "to"
$ws.Range("K29").Value = 0.006936416184971098
$ws.Range("L29").Value = 30
$ws.Range("M29").Value = 33
$ws.Range("N29").Value = 0.91
$ws.Range("O29").Value = 0.08999999999999997
$ws.Range("P29").Value = $true
$ws.Range("Q29").Value = 4295
$ws.Range("J30").Value = "the"
$ws.Range("K30").Value = 0.0052285050348567
$ws.Range("L30").Value = 27
$ws.Range("M30").Value = 28
$ws.Range("N30").Value = 0.96
$ws.Range("O30").Value = 0.04000000000000004
$ws.Range("P30").Value = $true
$ws.Range("Q30").Value = 5137
